# Re-generate the quadratic/linear bilevel experiment values.
# (volver a generar problemas cuadraticos y lineales)
#
# Only the data values change (the follower-constraint expressions/evaluations,
# the modified point, and the vec_bf / vec_BF vectors); headers and sheet
# layout stay the same.

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower ---------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

# row 2: J_0_L0_v  (expression / evaluation / type / miu / lambda / beta / gamma)
$ws.Range("A2").Value = "8.600000000000001 - y_1"
$ws.Range("B2").Value = "-8.600000000000001"
$ws.Range("C2").Value = "J_0_L0_v"
$ws.Range("D2").Value = "0.75"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "1.1"

# row 3: J_0_L0_v
$ws.Range("A3").Value = "-8.600000000000001 + y_1"
$ws.Range("B3").Value = "4.600000000000001"
$ws.Range("C3").Value = "J_0_L0_v"
$ws.Range("D3").Value = "0.19"
$ws.Range("E3").Value = "7.1"
$ws.Range("F3").Value = "0"

# row 4: J_0_LP_v
$ws.Range("A4").Value = "-5.000000000000002 - 2x + y_1 + 4y_2"
$ws.Range("B4").Value = "-10.999999999999998"
$ws.Range("C4").Value = "J_0_LP_v"
$ws.Range("D4").Value = "0.24"
$ws.Range("E4").Value = "5.0"
$ws.Range("F4").Value = "0"

# row 5: J_Ne_L0_v
$ws.Range("A5").Value = "-65.87 + 8x + y_1"
$ws.Range("B5").Value = "17.4"
$ws.Range("C5").Value = "J_Ne_L0_v"
$ws.Range("D5").Value = "0.92"
$ws.Range("E5").Value = "6.1"
$ws.Range("F5").Value = "8.8"

# row 6: J_Ne_L0_v
$ws.Range("A6").Value = "-7.400000000000002 - 2x - 2y_1"
$ws.Range("B6").Value = "-19.400000000000002"
$ws.Range("C6").Value = "J_Ne_L0_v"
$ws.Range("D6").Value = "0.82"
$ws.Range("E6").Value = "-7.9"
$ws.Range("F6").Value = "-1.0"

# --- Punto_modificado ---------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2").Value = "7.1"
$ws.Range("B2").Value = "8.600000000000001"
$ws.Range("C2").Value = "2.65"

# --- Vector_bf ------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_bf")
$ws.Range("A2").Value = "2.04"
$ws.Range("A3").Value = "-0.96"

# --- Vector_BF --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_BF")
$ws.Range("A2").Value = "-53.599999999999994"
$ws.Range("A3").Value = "-31.0"
$ws.Range("A4").Value = "-22.0"
